$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.293.15'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '3.496.56'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''588.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = '''133.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''0.485'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").Value = '''7.68'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.80%  '
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("E11").Value = '  +1.79%  '
$ws.Range("D12").Value = '4.092.34'
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("D14").Value = '''0.0000179'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.61%  '
$ws.Range("D15").Value = '3.497.39'
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").Value = '64.239.44'
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '''24.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.84%  '
$ws.Range("D18").Value = '''10.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.64%  '
$ws.Range("D19").Value = '''5.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").Value = '''13.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.47%  '
$ws.Range("D21").Value = '''385.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("E22").Value = '  +2.13%  '
$ws.Range("D23").Value = '3.636.05'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '''74.31'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("D25").Value = '''0.999'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.12%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").Value = '  +1.62%  '
$ws.Range("D28").Value = '''0.997'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("D29").Value = '''7.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '''2.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").Value = '''1.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("E32").Value = '  -1.63%  '
$ws.Range("E33").Value = '  +3.94%  '
$ws.Range("D34").Value = '3.524.65'
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").Value = '''23.24'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.97%  '
$ws.Range("E37").Value = '  +3.89%  '
$ws.Range("D38").Value = '''6.88'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("E39").Value = '  -0.61%  '
$ws.Range("D40").Value = '''164.43'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.91%  '
$ws.Range("D41").Value = '''0.0784'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("D42").Value = '''0.806'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").Value = '''4.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("D46").Value = '''24.24'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.82%  '
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("D48").Value = '2.434.78'
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("D49").Value = '''6.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("D50").Value = '''0.919'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.03%  '
$ws.Range("D51").Value = '''0.0258'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.89%  '
